$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DOMA-3100 add formatter convert to number for some colomns
# Append ":formatN()" to the template placeholders in columns C..H (processing,
# completed, canceled, deferred, closed, new_or_reopened) for both ticket rows,
# and switch those cells' number format to an integer format ("0") so the
# exported numbers render as numbers instead of generic/text.

$ws.Range("C2").Value = "{d.tickets[i].processing:formatN()}"
$ws.Range("D2").Value = "{d.tickets[i].completed:formatN()}"
$ws.Range("E2").Value = "{d.tickets[i].canceled:formatN()}"
$ws.Range("F2").Value = "{d.tickets[i].deferred:formatN()}"
$ws.Range("G2").Value = "{d.tickets[i].closed:formatN()}"
$ws.Range("H2").Value = "{d.tickets[i].new_or_reopened:formatN()}"

$ws.Range("C3").Value = "{d.tickets[i+1].processing:formatN()}"
$ws.Range("D3").Value = "{d.tickets[i+1].completed:formatN()}"
$ws.Range("E3").Value = "{d.tickets[i+1].canceled:formatN()}"
$ws.Range("F3").Value = "{d.tickets[i+1].deferred:formatN()}"
$ws.Range("G3").Value = "{d.tickets[i+1].closed:formatN()}"
$ws.Range("H3").Value = "{d.tickets[i+1].new_or_reopened:formatN()}"

$ws.Range("C2:H2").NumberFormat = "0"
$ws.Range("C3:H3").NumberFormat = "0"
